$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.771.10'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.623.74'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.74'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2557'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06301'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.34'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07777'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').Value = '1.627.43'
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').Value = '1.846.63'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5510'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.36'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '0.0₅7477'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '25.771.63'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.400'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '193.66'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.766'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.986'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.80%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.878'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '141.60'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1251'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.42%  '
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.704'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.236'
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.04862'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.226'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.148'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.535'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.369'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.8931'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.536'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.93%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5496'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('D39').Value = '1.110.97'
$ws.Range('E39').Value = '  -2.78%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01543'
$ws.Range('D40').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.535'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.7957'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.61%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '97.21'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.92%  '
$ws.Range('D45').Value = '1.771.96'
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('D46').Value = '0.0₈116'
$ws.Range('E46').Value = '  -8.77%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4424'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.30%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9990'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '54.53'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05122'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.543'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.11%  '
